$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.988.02'
$ws.Range('E2').Value = '  +1.10%  '

$ws.Range('D3').Value = '2.641.63'
$ws.Range('E3').Value = '  +1.59%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = '530.85'
$ws.Range('E5').Value = '  +4.04%  '

$ws.Range('D6').Value = '155.51'
$ws.Range('E6').Value = '  +0.68%  '

$ws.Range('E7').Value = '  +0.04%  '

$ws.Range('E8').Value = '  +0.71%  '

$ws.Range('E9').Value = '  -0.27%  '

$ws.Range('E10').Value = '  +5.02%  '

$ws.Range('D11').Value = '0.351'
$ws.Range('E11').Value = '  +1.68%  '

$ws.Range('E12').Value = '  -0.11%  '

$ws.Range('D13').Value = '3.106.92'
$ws.Range('E13').Value = '  +1.69%  '

$ws.Range('D14').Value = '61.040.32'

$ws.Range('D15').Value = '22.03'
$ws.Range('E15').Value = '  +2.34%  '

$ws.Range('E16').Value = '  +2.88%  '

$ws.Range('D17').Value = '2.642.03'
$ws.Range('E17').Value = '  +1.43%  '

$ws.Range('D18').Value = '4.76'
$ws.Range('E18').Value = '  +0.44%  '

$ws.Range('D19').Value = '353.86'
$ws.Range('E19').Value = '  +0.62%  '

$ws.Range('D20').Value = '10.66'
$ws.Range('E20').Value = '  +0.88%  '

$ws.Range('D21').Value = '6.23'
$ws.Range('E21').Value = '  +1.65%  '

$ws.Range('E22').Value = '  +0.39%  '

$ws.Range('D23').Value = '61.77'
$ws.Range('E23').Value = '  +2.21%  '

$ws.Range('D24').Value = '0.430'
$ws.Range('E24').Value = '  +2.40%  '

$ws.Range('D25').Value = '0.168'
$ws.Range('E25').Value = '  +1.67%  '

$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  +0.12%  '

$ws.Range('D27').Value = '0.0₃0863'
$ws.Range('E27').Value = '  +2.69%  '

$ws.Range('D28').Value = '7.39'
$ws.Range('E28').Value = '  +0.37%  '

$ws.Range('E29').Value = '  -0.04%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '1.63'
$ws.Range('E30').Value = '  +4.30%  '

$ws.Range('D31').Value = '19.51'
$ws.Range('E31').Value = '  +0.72%  '

$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D32').Value = '6.13'
$ws.Range('E32').Value = '  +6.98%  '

$ws.Range('D33').Value = '150.51'
$ws.Range('E33').Value = '  -0.46%  '

$ws.Range('D34').Value = '4.14'
$ws.Range('E34').Value = '  +3.71%  '

$ws.Range('D35').Value = '1.20'
$ws.Range('E35').Value = '  +1.80%  '

$ws.Range('D36').Value = '0.920'
$ws.Range('E36').Value = '  +9.17%  '

$ws.Range('E37').Value = '  +1.53%  '

$ws.Range('D38').Value = '307.80'
$ws.Range('E38').Value = '  +4.82%  '

$ws.Range('E39').Value = '  +1.24%  '

$ws.Range('E40').Value = '  +1.62%  '

$ws.Range('D41').Value = '0.645'
$ws.Range('E41').Value = '  +3.75%  '

$ws.Range('E42').Value = '  +1.57%  '

$ws.Range('D43').Value = '0.0562'
$ws.Range('E43').Value = '  +1.74%  '

$ws.Range('D44').Value = '0.998'
$ws.Range('E44').Value = '  +0.05%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '19.86'
$ws.Range('E45').Value = '  +0.94%  '

$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = '4.96'
$ws.Range('E46').Value = '  +1.68%  '

$ws.Range('E47').Value = '  +2.36%  '

$ws.Range('D48').Value = '19.26'
$ws.Range('E48').Value = '  +7.82%  '

$ws.Range('E49').Value = '  +0.22%  '

$ws.Range('D50').Value = '1.987.69'
$ws.Range('E50').Value = '  -0.49%  '
